# Weekly CompStat (108th Precinct) refresh: new crime data collected.
# Updates the report header (week number + date range) and the weekly
# crime-statistics grid (Week to Date / 28 Day / Year to Date / 2 Year columns).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Report header: "Volume 32  Number 51" -> "...Number 52",
# and the covered week "12/15/2025 ... 12/21/2025" -> "12/22/2025 ... 12/28/2025" ---
$ws.Range("A8").Value = "Volume 32   Number  52"
$ws.Range("C9").Value = "Report Covering the Week  12/22/2025  Through  12/28/2025"

# Row 15 - Rape
$ws.Range("F15").Value = 1
$ws.Range("L15").Value = -12

# Row 16 - Robbery
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 300
$ws.Range("F16").Value = 16
$ws.Range("G16").Value = 11
$ws.Range("H16").Value = 45.454545454545
$ws.Range("I16").Value = 213
$ws.Range("J16").Value = 232
$ws.Range("K16").Value = -8.189655172413
$ws.Range("L16").Value = -19.011406844106
$ws.Range("M16").Value = 10.362694300518
$ws.Range("N16").Value = -79.056047197640

# Row 17 - Fel. Assault
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 18
$ws.Range("G17").Value = 15
$ws.Range("H17").Value = 20
$ws.Range("I17").Value = 297
$ws.Range("J17").Value = 264
$ws.Range("K17").Value = 12.5
$ws.Range("L17").Value = 25.847457627118
$ws.Range("M17").Value = 128.461538461538
$ws.Range("N17").Value = -5.414012738853

# Row 18 - Burglary
$ws.Range("C18").Value = 6
$ws.Range("D18").Value = 10
$ws.Range("E18").Value = -40
$ws.Range("F18").Value = 16
$ws.Range("G18").Value = 24
$ws.Range("H18").Value = -33.333333333333
$ws.Range("I18").Value = 231
$ws.Range("J18").Value = 250
$ws.Range("K18").Value = -7.6
$ws.Range("L18").Value = 0.873362445414
$ws.Range("M18").Value = -9.411764705882
$ws.Range("N18").Value = -85.087153001936

# Row 19 - Gr. Larceny
$ws.Range("C19").Value = 10
$ws.Range("E19").Value = -28.571428571428
$ws.Range("F19").Value = 55
$ws.Range("G19").Value = 56
$ws.Range("H19").Value = -1.785714285714
$ws.Range("I19").Value = 796
$ws.Range("J19").Value = 735
$ws.Range("K19").Value = 8.299319727891
$ws.Range("L19").Value = 6.133333333333
$ws.Range("M19").Value = 74.945054945054
$ws.Range("N19").Value = -12.527472527472

# Row 20 - G.L.A.
$ws.Range("C20").Value = 6
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 200
$ws.Range("F20").Value = 27
$ws.Range("G20").Value = 15
$ws.Range("H20").Value = 80
$ws.Range("I20").Value = 253
$ws.Range("J20").Value = 255
$ws.Range("K20").Value = -0.784313725490
$ws.Range("L20").Value = -18.910256410256
$ws.Range("M20").Value = 30.412371134020
$ws.Range("N20").Value = -88.015158692562

# Row 21 - TOTAL (7 major crimes)
$ws.Range("D21").Value = 29
$ws.Range("E21").Value = -3.448275862068
$ws.Range("F21").Value = 133
$ws.Range("G21").Value = 121
$ws.Range("H21").Value = 9.917355371900
$ws.Range("I21").Value = 1814
$ws.Range("J21").Value = 1758
$ws.Range("K21").Value = 3.185437997724
$ws.Range("L21").Value = -0.055096418732
$ws.Range("M21").Value = 45.585874799357
$ws.Range("N21").Value = -69.373628228938

# Row 22 - Transit
$ws.Range("C22").Value = 2
$ws.Range("C22").NumberFormat = "#,##0"  # was a blank-dash placeholder; now a real number
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 100
$ws.Range("F22").Value = 5
$ws.Range("H22").Value = -28.571428571428
$ws.Range("I22").Value = 66
$ws.Range("J22").Value = 63
$ws.Range("K22").Value = 4.761904761904
$ws.Range("L22").Value = -23.255813953488
$ws.Range("M22").Value = 34.693877551020

# Row 24 - Petit Larceny
$ws.Range("C24").Value = 19
$ws.Range("D24").Value = 30
$ws.Range("E24").Value = -36.666666666666
$ws.Range("F24").Value = 97
$ws.Range("G24").Value = 166
$ws.Range("H24").Value = -41.566265060241
$ws.Range("I24").Value = 1540
$ws.Range("J24").Value = 2055
$ws.Range("K24").Value = -25.060827250608
$ws.Range("L24").Value = -25.747348119575
$ws.Range("M24").Value = 65.948275862069

# Row 25 - Retail Theft
$ws.Range("C25").Value = 14
$ws.Range("D25").Value = 13
$ws.Range("E25").Value = 7.692307692307
$ws.Range("F25").Value = 52
$ws.Range("G25").Value = 101
$ws.Range("H25").Value = -48.514851485148
$ws.Range("I25").Value = 795
$ws.Range("J25").Value = 1379
$ws.Range("K25").Value = -42.349528643944
$ws.Range("L25").Value = -38.084112149532

# Row 26 - Misd. Assault
$ws.Range("C26").Value = 9
$ws.Range("D26").Value = 4
$ws.Range("E26").Value = 125
$ws.Range("F26").Value = 42
$ws.Range("G26").Value = 37
$ws.Range("H26").Value = 13.513513513513
$ws.Range("I26").Value = 532
$ws.Range("J26").Value = 551
$ws.Range("K26").Value = -3.448275862068
$ws.Range("L26").Value = 1.720841300191
$ws.Range("M26").Value = 7.042253521126

# Row 27 - UCR Rape*
$ws.Range("F27").Value = 1
$ws.Range("L27").Value = -12.121212121212

# Row 28 - Other Sex Crimes
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 0
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = -20
$ws.Range("I28").Value = 69
$ws.Range("J28").Value = 65
$ws.Range("K28").Value = 6.153846153846
$ws.Range("L28").Value = -18.823529411764

# Row 31 - Hate Crimes
$ws.Range("D31").Value = 1
$ws.Range("D31").NumberFormat = "#,##0"  # was a blank-dash placeholder; now a real number
$ws.Range("E31").Value = -100
$ws.Range("E31").NumberFormat = "#,##0.0;""-""#,##0.0"  # was a blank-dash placeholder; now a real number
$ws.Range("G31").Value = 2
$ws.Range("H31").Value = -50
$ws.Range("J31").Value = 11
$ws.Range("K31").Value = -9.090909090909

# Row 33 - Traffic Fatalities
$ws.Range("L33").Value = -85.714285714285
